$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = 'https://www.alojamiento.io/property/apartamentos-centro-col%c3%b3n/BC-189483'
$ws.Cells.Item(2,2).Value = 'H1 Tag Existence'
$ws.Cells.Item(2,3).Value = $true
$ws.Cells.Item(2,4).Value = 'H1 tag found'
$ws.Cells.Item(2,5).Value = '2024-12-06 17:44:12'
$ws.Cells.Item(2,6).Value = ''
$ws.Cells.Item(2,7).Value = $true

# Row 3
$ws.Cells.Item(3,1).Value = 'https://www.alojamiento.io/property/apartamentos-centro-col%c3%b3n/BC-189483'
$ws.Cells.Item(3,2).Value = 'H1 Tag Existence'
$ws.Cells.Item(3,3).Value = $true
$ws.Cells.Item(3,4).Value = 'H1 tag found'
$ws.Cells.Item(3,5).Value = '2024-12-06 17:58:08'
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = $false

# Row 4
$ws.Cells.Item(4,1).Value = 'https://www.alojamiento.io/property/apartamentos-centro-col%c3%b3n/BC-189483'
$ws.Cells.Item(4,2).Value = 'H1-H6 Tag Sequence'
$ws.Cells.Item(4,3).Value = $false
$ws.Cells.Item(4,4).Value = 'H5 missing'
$ws.Cells.Item(4,5).Value = '2024-12-06 17:58:08'
$ws.Cells.Item(4,6).Value = ''
$ws.Cells.Item(4,7).Value = $true

# Row 5
$ws.Cells.Item(5,1).Value = 'https://www.alojamiento.io/property/apartamentos-centro-col%c3%b3n/BC-189483'
$ws.Cells.Item(5,2).Value = 'H1 Tag Existence'
$ws.Cells.Item(5,3).Value = $true
$ws.Cells.Item(5,4).Value = 'H1 tag found'
$ws.Cells.Item(5,5).Value = '2024-12-06 18:01:35'
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = $false

# Row 6
$ws.Cells.Item(6,1).Value = 'https://www.alojamiento.io/property/apartamentos-centro-col%c3%b3n/BC-189483'
$ws.Cells.Item(6,2).Value = 'H1-H6 Tag Sequence'
$ws.Cells.Item(6,3).Value = $false
$ws.Cells.Item(6,4).Value = 'H5 missing'
$ws.Cells.Item(6,5).Value = '2024-12-06 18:01:35'
$ws.Cells.Item(6,6).Value = 0
$ws.Cells.Item(6,7).Value = $false

# Row 7
$ws.Cells.Item(7,1).Value = 'https://www.alojamiento.io/property/apartamentos-centro-col%c3%b3n/BC-189483'
$ws.Cells.Item(7,2).Value = 'H1 Tag Existence'
$ws.Cells.Item(7,3).Value = $true
$ws.Cells.Item(7,4).Value = 'H1 tag found'
$ws.Cells.Item(7,5).Value = '2024-12-06 18:02:50'
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = $false

# Row 8
$ws.Cells.Item(8,1).Value = 'https://www.alojamiento.io/property/apartamentos-centro-col%c3%b3n/BC-189483'
$ws.Cells.Item(8,2).Value = 'H1-H6 Tag Sequence'
$ws.Cells.Item(8,3).Value = $false
$ws.Cells.Item(8,4).Value = 'H5 missing'
$ws.Cells.Item(8,5).Value = '2024-12-06 18:02:50'
$ws.Cells.Item(8,6).Value = 0
$ws.Cells.Item(8,7).Value = $false

# Row 9
$ws.Cells.Item(9,1).Value = 'https://www.alojamiento.io/property/apartamentos-centro-col%c3%b3n/BC-189483'
$ws.Cells.Item(9,2).Value = 'H1 Tag Existence'
$ws.Cells.Item(9,3).Value = $true
$ws.Cells.Item(9,4).Value = 'H1 tag found'
$ws.Cells.Item(9,5).Value = '2024-12-06 18:03:52'
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = $false

# Row 10
$ws.Cells.Item(10,1).Value = 'https://www.alojamiento.io/property/apartamentos-centro-col%c3%b3n/BC-189483'
$ws.Cells.Item(10,2).Value = 'H1-H6 Tag Sequence'
$ws.Cells.Item(10,3).Value = $false
$ws.Cells.Item(10,4).Value = 'H5 missing'
$ws.Cells.Item(10,5).Value = '2024-12-06 18:03:52'
$ws.Cells.Item(10,6).Value = 0
$ws.Cells.Item(10,7).Value = $false

# Row 11
$ws.Cells.Item(11,1).Value = 'https://www.alojamiento.io/property/apartamentos-centro-col%c3%b3n/BC-189483'
$ws.Cells.Item(11,2).Value = 'H1 Tag Existence'
$ws.Cells.Item(11,3).Value = $true
$ws.Cells.Item(11,4).Value = 'H1 tag found'
$ws.Cells.Item(11,5).Value = '2024-12-06 18:05:25'
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = $false

# Row 12
$ws.Cells.Item(12,1).Value = 'https://www.alojamiento.io/property/apartamentos-centro-col%c3%b3n/BC-189483'
$ws.Cells.Item(12,2).Value = 'H1-H6 Tag Sequence'
$ws.Cells.Item(12,3).Value = $false
$ws.Cells.Item(12,4).Value = 'H5 missing'
$ws.Cells.Item(12,5).Value = '2024-12-06 18:05:25'
$ws.Cells.Item(12,6).Value = 0
$ws.Cells.Item(12,7).Value = $false

# Row 13
$ws.Cells.Item(13,1).Value = 'https://www.alojamiento.io/property/apartamentos-centro-col%c3%b3n/BC-189483'
$ws.Cells.Item(13,2).Value = 'H1 Tag Existence'
$ws.Cells.Item(13,3).Value = $true
$ws.Cells.Item(13,4).Value = 'H1 tag found'
$ws.Cells.Item(13,5).Value = '2024-12-06 18:08:48'
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = $false

# Row 14
$ws.Cells.Item(14,1).Value = 'https://www.alojamiento.io/property/apartamentos-centro-col%c3%b3n/BC-189483'
$ws.Cells.Item(14,2).Value = 'H1-H6 Tag Sequence'
$ws.Cells.Item(14,3).Value = $true
$ws.Cells.Item(14,4).Value = 'Tags found in correct sequence'
$ws.Cells.Item(14,5).Value = '2024-12-06 18:08:48'
$ws.Cells.Item(14,6).Value = 0
$ws.Cells.Item(14,7).Value = $true

# Row 15
$ws.Cells.Item(15,1).Value = 'https://www.alojamiento.io/property/consultar-disponibilidad/BC-4505653'
$ws.Cells.Item(15,2).Value = 'H1 Tag Existence'
$ws.Cells.Item(15,3).Value = $true
$ws.Cells.Item(15,4).Value = 'H1 tag found'
$ws.Cells.Item(15,5).Value = '2024-12-06 18:11:51'
$ws.Cells.Item(15,6).Value = ''
$ws.Cells.Item(15,7).Value = $true

# Row 16
$ws.Cells.Item(16,1).Value = 'https://www.alojamiento.io/property/consultar-disponibilidad/BC-4505653'
$ws.Cells.Item(16,2).Value = 'H1-H6 Tag Sequence'
$ws.Cells.Item(16,3).Value = $true
$ws.Cells.Item(16,4).Value = 'Tags found in correct sequence'
$ws.Cells.Item(16,5).Value = '2024-12-06 18:11:51'
$ws.Cells.Item(16,6).Value = ''
$ws.Cells.Item(16,7).Value = $true

# Row 17
$ws.Cells.Item(17,1).Value = 'https://www.alojamiento.io/property/consultar-disponibilidad/BC-4505653'
$ws.Cells.Item(17,2).Value = 'H1 Tag Existence'
$ws.Cells.Item(17,3).Value = $true
$ws.Cells.Item(17,4).Value = 'H1 tag found'
$ws.Cells.Item(17,5).Value = '2024-12-06 18:13:42'
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = $false

# Row 18
$ws.Cells.Item(18,1).Value = 'https://www.alojamiento.io/property/consultar-disponibilidad/BC-4505653'
$ws.Cells.Item(18,2).Value = 'H1-H6 Tag Sequence'
$ws.Cells.Item(18,3).Value = $true
$ws.Cells.Item(18,4).Value = 'Tags found in correct sequence'
$ws.Cells.Item(18,5).Value = '2024-12-06 18:13:42'
$ws.Cells.Item(18,6).Value = 1
$ws.Cells.Item(18,7).Value = $false

# Row 19
$ws.Cells.Item(19,1).Value = 'https://www.alojamiento.io/property/bonita-casa-de-campo-t%C3%ADpica-mallorquina/BC-12224317'
$ws.Cells.Item(19,2).Value = 'H1 Tag Existence'
$ws.Cells.Item(19,3).Value = $true
$ws.Cells.Item(19,4).Value = 'H1 tag found'
$ws.Cells.Item(19,5).Value = '2024-12-06 18:17:02'
$ws.Cells.Item(19,6).Value = ''
$ws.Cells.Item(19,7).Value = $true

# Row 20
$ws.Cells.Item(20,1).Value = 'https://www.alojamiento.io/property/bonita-casa-de-campo-t%C3%ADpica-mallorquina/BC-12224317'
$ws.Cells.Item(20,2).Value = 'H1-H6 Tag Sequence'
$ws.Cells.Item(20,3).Value = $true
$ws.Cells.Item(20,4).Value = 'Tags found in correct sequence'
$ws.Cells.Item(20,5).Value = '2024-12-06 18:17:02'
$ws.Cells.Item(20,6).Value = ''
$ws.Cells.Item(20,7).Value = $true

# Row 21
$ws.Cells.Item(21,1).Value = 'https://www.alojamiento.io/property/bonita-casa-de-campo-t%C3%ADpica-mallorquina/BC-12224317'
$ws.Cells.Item(21,2).Value = 'H1 Tag Existence'
$ws.Cells.Item(21,3).Value = $true
$ws.Cells.Item(21,4).Value = 'H1 tag found'
$ws.Cells.Item(21,5).Value = '2024-12-06 18:25:50'
$ws.Cells.Item(21,6).Value = 1
$ws.Cells.Item(21,7).Value = $false

# Row 22
$ws.Cells.Item(22,1).Value = 'https://www.alojamiento.io/'
$ws.Cells.Item(22,2).Value = 'H1 Tag Existence'
$ws.Cells.Item(22,3).Value = $true
$ws.Cells.Item(22,4).Value = 'H1 tag found'
$ws.Cells.Item(22,5).Value = '2024-12-06 18:29:47'
$ws.Cells.Item(22,6).Value = ''
$ws.Cells.Item(22,7).Value = $true

# Row 23
$ws.Cells.Item(23,1).Value = 'https://www.alojamiento.io/'
$ws.Cells.Item(23,2).Value = 'H1-H6 Tag Sequence'
$ws.Cells.Item(23,3).Value = $true
$ws.Cells.Item(23,4).Value = 'Tags found in correct sequence'
$ws.Cells.Item(23,5).Value = '2024-12-06 18:29:47'
$ws.Cells.Item(23,6).Value = ''
$ws.Cells.Item(23,7).Value = $true

# Row 24
$ws.Cells.Item(24,1).Value = 'https://www.alojamiento.io/'
$ws.Cells.Item(24,2).Value = 'H1 Tag Existence'
$ws.Cells.Item(24,3).Value = $true
$ws.Cells.Item(24,4).Value = 'H1 tag found'
$ws.Cells.Item(24,5).Value = '2024-12-06 18:32:30'
$ws.Cells.Item(24,6).Value = $true
$ws.Cells.Item(24,7).Value = $false

# Row 25
$ws.Cells.Item(25,1).Value = 'https://www.alojamiento.io/'
$ws.Cells.Item(25,2).Value = 'HTML Tag Sequence'
$ws.Cells.Item(25,3).Value = $false
$ws.Cells.Item(25,4).Value = 'Improper heading hierarchy'
$ws.Cells.Item(25,5).Value = '2024-12-06 18:32:30'
$ws.Cells.Item(25,6).Value = 'N/A'
$ws.Cells.Item(25,7).Value = $true
